$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 295
$ws1.Range("F3").Value = 1020
$ws1.Range("F6").Value = 3255
$ws1.Range("F7").Value = 51
$ws1.Range("F13").Value = 43
$ws1.Range("F14").Value = 103
$ws1.Range("F16").Value = 1508
$ws1.Range("F17").Value = 1508
$ws1.Range("F18").Value = 11
$ws1.Range("F19").Value = 280
$ws1.Range("F20").Value = 24
$ws1.Range("G20").Value = 66
$ws1.Range("F21").Value = 593
$ws1.Range("F22").Value = 318
$ws1.Range("F24").Value = 542
$ws1.Range("F25").Value = 43265
$ws1.Range("F26").Value = 43267
$ws1.Range("F27").Value = 708
$ws1.Range("F29").Value = 32219
$ws1.Range("F30").Value = 32219
$ws1.Range("F31").Value = 426
$ws1.Range("F35").Value = 916
$ws1.Range("F36").Value = 221
$ws1.Range("F38").Value = 480
$ws1.Range("F39").Value = 1139
$ws1.Range("F40").Value = 5289
$ws1.Range("F41").Value = 680
$ws1.Range("F42").Value = 407
$ws1.Range("F45").Value = 309
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 213
$ws2.Range("F15").Value = 741
$ws2.Range("F19").Value = 5
$ws2.Range("F23").Value = 459
$ws2.Range("F35").Value = 1130
$ws2.Range("F45").Value = 794
$ws2.Range("F47").Value = 56
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 688
$ws3.Range("F5").Value = 519
$ws3.Range("F6").Value = 515
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 688
$ws4.Range("F3").Value = 295
$ws4.Range("F4").Value = 519
$ws4.Range("F6").Value = 1020
$ws4.Range("F8").Value = 213
$ws4.Range("F10").Value = 3255
$ws4.Range("F11").Value = 51
$ws4.Range("F15").Value = 515
$ws4.Range("F20").Value = 43
$ws4.Range("F21").Value = 103
$ws4.Range("F23").Value = 1508
$ws4.Range("F24").Value = 1508
$ws4.Range("F25").Value = 280
$ws4.Range("F27").Value = 24
$ws4.Range("G27").Value = 66
$ws4.Range("F28").Value = 593
$ws4.Range("F30").Value = 318
$ws4.Range("F31").Value = 542
$ws4.Range("F32").Value = 43267
$ws4.Range("F34").Value = 708
$ws4.Range("F36").Value = 32219
$ws4.Range("F37").Value = 426
$ws4.Range("F38").Value = 916
$ws4.Range("F39").Value = 221
$ws4.Range("F41").Value = 480
$ws4.Range("F42").Value = 1139
$ws4.Range("F43").Value = 5289
$ws4.Range("F44").Value = 680
$ws4.Range("F46").Value = 407
$ws4.Range("F49").Value = 309
$ws4.Range("F51").Value = 794
$ws4.Range("F53").Value = 56
